$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "6"
$ws.Range("D2").Value = 0.01485
$ws.Range("E2").Value = -0.062
$ws.Range("F2").Value = 0.0216
$ws.Range("I2").Value = 0.000001660580936412585
$ws.Range("J2").Value = 0.000001334264041243421
$ws.Range("K2").Value = 2335.1
$ws.Range("L2").Value = 0.1494703152504401
$ws.Range("M2").Value = 1354.339
$ws.Range("N2").Value = 0.03718545792011246
$ws.Range("O2").Value = 0.5799918633034987
$ws.Range("P2").Value = 1169.639
$ws.Range("Q2").Value = 0.0321142356649424
$ws.Range("R2").Value = 0.50089460836795
$ws.Range("S2").Value = 184.7
$ws.Range("T2").Value = 0.1363764906718333
$ws.Range("U2").Value = 11723.4
$ws.Range("V2").Value = 0.3218839576949689
$ws.Range("W2").Value = 0.06233712740758868
$ws.Range("X2").Value = 0.09416848029118631
$ws.Range("Y2").Value = -0.03183135288359763
$ws.Range("Z2").Value = 0.2653619908478321
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0.05761080692926477
$ws.Range("AC2").Value = -0.05761080692926477
$ws.Range("AD2").Value = 36910.8
$ws.Range("AE2").Value = 0.03028787160447198
$ws.Range("AF2").Value = 36910.83028787161
$ws.Range("AG2").Value = 25187.4302878716
$ws.Range("AH2").Value = 0.5033384476466117
$ws.Range("AI2").Value = 0.5341628633915537
$ws.Range("AJ2").Value = 0.4088295774501264
$ws.Range("AK2").Value = 0.4389818374998662
$ws.Range("AN2").Value = 1153462.5
$ws.Range("AP2").Value = 787107.1964959876
# Row 3
$ws.Range("B3").Value = "Capitec Bank Holdings Limited (JSE:CPI)"
$ws.Range("D3").Value = 0.102
$ws.Range("E3").Value = 0.0653
$ws.Range("F3").Value = 0.122
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 231.1
$ws.Range("L3").Value = 0.2867601439384539
$ws.Range("M3").Value = 13.4
$ws.Range("N3").Value = 0.001185389630494591
$ws.Range("O3").Value = 0.05798355690177413
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 13.4
$ws.Range("T3").Value = 1.0
$ws.Range("U3").Value = 3054.1
$ws.Range("V3").Value = 0.2701715276487708
$ws.Range("W3").Value = 0.1521696187528807
$ws.Range("X3").Value = 0.04660605020090043
$ws.Range("Y3").Value = 0.1055635685519803
$ws.Range("Z3").Value = -1.683975176045301
$ws.Range("AA3").Value = 0
$ws.Range("AB3").Value = 0.04636309190295776
$ws.Range("AC3").Value = -0.04636309190295776
$ws.Range("AD3").Value = 316.7
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 316.7
$ws.Range("AG3").Value = -2737.4
$ws.Range("AH3").Value = 0.02725238791842354
$ws.Range("AI3").Value = 0.1708198489751888
$ws.Range("AJ3").Value = -0.3195321528207403
$ws.Range("AK3").Value = 2.280976585284559
$ws.Range("AN3").Value = ""
$ws.Range("AP3").Value = ""
# Row 4
$ws.Range("B4").Value = "Standard Bank Group Limited (JSE:SBK)"
$ws.Range("D4").Value = 0.0222
$ws.Range("E4").Value = -0.062
$ws.Range("F4").Value = 0.0274
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 973.8
$ws.Range("L4").Value = 0.1493420850841947
$ws.Range("M4").Value = 2.36
$ws.Range("N4").Value = 0.0001710293648723077
$ws.Range("O4").Value = 0.002423495584308893
$ws.Range("P4").Value = 2.36
$ws.Range("Q4").Value = 0.0001710293648723077
$ws.Range("R4").Value = 0.002423495584308893
$ws.Range("U4").Value = 6190.8
$ws.Range("V4").Value = 0.448647708496391
$ws.Range("W4").Value = 0.07683810185112125
$ws.Range("X4").Value = 0.06320525577116365
$ws.Range("Y4").Value = 0.0136328460799576
$ws.Range("Z4").Value = 0.4622768585080892
$ws.Range("AA4").Value = 0
$ws.Range("AB4").Value = 0.0530752498901797
$ws.Range("AC4").Value = -0.0530752498901797
$ws.Range("AD4").Value = 9086.0
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 9086.0
$ws.Range("AG4").Value = 2895.2
$ws.Range("AH4").Value = 0.3970320911696847
$ws.Range("AI4").Value = 0.4191442740168377
$ws.Range("AJ4").Value = 0.1734275787708159
$ws.Range("AK4").Value = 0.1869475097987305
$ws.Range("AN4").Value = ""
$ws.Range("AP4").Value = ""
# Row 5
$ws.Range("D5").Value = 0.0131
$ws.Range("E5").Value = -0.125
$ws.Range("F5").Value = -0.00411
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 409.2
$ws.Range("L5").Value = 0.1135373602286285
$ws.Range("M5").Value = 648.3
$ws.Range("N5").Value = 0.09569568683021874
$ws.Range("O5").Value = 1.584310850439883
$ws.Range("P5").Value = 537.9
$ws.Range("Q5").Value = 0.07939952174298112
$ws.Range("R5").Value = 1.314516129032258
$ws.Range("S5").Value = 110.4
$ws.Range("T5").Value = 0.1702915316982878
$ws.Range("U5").Value = 957.6
$ws.Range("V5").Value = 0.1413515189088655
$ws.Range("W5").Value = 0.05138186064616582
$ws.Range("X5").Value = 0.09199688565210991
$ws.Range("Y5").Value = -0.04061502500594408
$ws.Range("Z5").Value = 0.1786711085332421
$ws.Range("AA5").Value = 0
$ws.Range("AB5").Value = 0.05742430414328496
$ws.Range("AC5").Value = -0.05742430414328496
$ws.Range("AD5").Value = 11869.0
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 11869.0
$ws.Range("AG5").Value = 10911.4
$ws.Range("AH5").Value = 0.6366259735244266
$ws.Range("AI5").Value = 0.6094699167620916
$ws.Range("AJ5").Value = 0.6169512608843153
$ws.Range("AK5").Value = 0.5892734666544254
$ws.Range("AN5").Value = ""
$ws.Range("AP5").Value = ""
# Row 6
$ws.Range("D6").Value = 0.0166
$ws.Range("E6").Value = -0.0917
$ws.Range("F6").Value = 0.0158
$ws.Range("K6").Value = 375.1
$ws.Range("L6").Value = 0.1491807190582246
$ws.Range("M6").Value = 625.0
$ws.Range("N6").Value = 0.1454909446436054
$ws.Range("O6").Value = 1.666222340709144
$ws.Range("P6").Value = 609.3
$ws.Range("Q6").Value = 0.141836212114158
$ws.Range("R6").Value = 1.62436683551053
$ws.Range("S6").Value = 15.70000000000005
$ws.Range("T6").Value = 0.02512000000000007
$ws.Range("U6").Value = 920.4
$ws.Range("V6").Value = 0.214255784719959
$ws.Range("W6").Value = 0.06188033060032665
$ws.Range("X6").Value = 0.09634007493026271
$ws.Range("Y6").Value = -0.03445974432993607
$ws.Range("Z6").Value = 0.1884956481974317
$ws.Range("AB6").Value = 0.05779730971524458
$ws.Range("AC6").Value = -0.05779730971524458
$ws.Range("AD6").Value = 8234.8
$ws.Range("AF6").Value = 8234.8
$ws.Range("AG6").Value = 7314.4
$ws.Range("AH6").Value = 0.6571752350246597
$ws.Range("AI6").Value = 0.5930887458047045
$ws.Range("AJ6").Value = 0.6299977605898262
$ws.Range("AK6").Value = 0.5641998734977861
# Row 7
$ws.Range("B7").Value = "Sasfin Holdings Limited (JSE:SFN)"
$ws.Range("D7").Value = -0.00658
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = ""
$ws.Range("I7").Value = 0.0005037364209535068
$ws.Range("J7").Value = 0.0005037364209535068
$ws.Range("K7").Value = -2.6
$ws.Range("L7").Value = -0.05048543689320389
$ws.Range("M7").Value = 2.7
$ws.Range("N7").Value = 0.05252918287937744
$ws.Range("O7").Value = -1.038461538461539
$ws.Range("P7").Value = 2.7
$ws.Range("Q7").Value = 0.05252918287937744
$ws.Range("R7").Value = -1.038461538461539
$ws.Range("S7").Value = 0
$ws.Range("T7").Value = 0
$ws.Range("U7").Value = 91.7
$ws.Range("V7").Value = 1.784046692607004
$ws.Range("W7").Value = -0.0231729055258467
$ws.Range("X7").Value = 0.2038594702319572
$ws.Range("Y7").Value = -0.2270323757578039
$ws.Range("Z7").Value = 0.11678526550147
$ws.Range("AA7").Value = 0.00005882899166381556
$ws.Range("AB7").Value = 0.06400698860822952
$ws.Range("AC7").Value = -0.0639481596165657
$ws.Range("AD7").Value = 308.4
$ws.Range("AE7").Value = 0.03028787160447198
$ws.Range("AF7").Value = 308.4302878716044
$ws.Range("AG7").Value = 216.7302878716044
$ws.Range("AH7").Value = 0.8571548818082244
$ws.Range("AI7").Value = 0.7564566505020841
$ws.Range("AJ7").Value = 0.8083021488992951
$ws.Range("AK7").Value = 0.6857896100124959
$ws.Range("AN7").Value = 9637.499999999998
$ws.Range("AP7").Value = 6772.821495987639
# Row 8
$ws.Range("A8").Value = "South Africa"
$ws.Range("B8").Value = "Nedbank Limited (JSE:NBKP)"
$ws.Range("C8").Value = "Bank (Money Center)"
$ws.Range("D8").Value = 0.0129
$ws.Range("E8").Value = -0.0575
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 348.5
$ws.Range("L8").Value = 0.1639228598306679
$ws.Range("M8").Value = 62.579
$ws.Range("N8").Value = 0.318792664289353
$ws.Range("O8").Value = 0.1795667144906743
$ws.Range("P8").Value = 17.379
$ws.Range("Q8").Value = 0.0885328578706062
$ws.Range("R8").Value = 0.04986800573888091
$ws.Range("S8").Value = 45.2
$ws.Range("T8").Value = 0.7222870291950975
$ws.Range("U8").Value = 508.8
$ws.Range("V8").Value = 2.591951095262353
$ws.Range("W8").Value = 0.06279392421485072
$ws.Range("X8").Value = 0.9976247840984154
$ws.Range("Y8").Value = -0.9348308598835646
$ws.Range("Z8").Value = 0.1882482113763548
$ws.Range("AA8").Value = 0
$ws.Range("AB8").Value = 0.07581836828137212
$ws.Range("AC8").Value = -0.07581836828137212
$ws.Range("AD8").Value = 7095.9
$ws.Range("AE8").Value = 0
$ws.Range("AF8").Value = 7095.9
$ws.Range("AG8").Value = 6587.099999999999
$ws.Range("AH8").Value = 0.9730808260881489
$ws.Range("AI8").Value = 0.6012353628984426
$ws.Range("AJ8").Value = 0.9710617094672288
$ws.Range("AK8").Value = 0.5832698744399384
$ws.Range("AL8").Value = 0
$ws.Range("AM8").Value = 0
